# Add a new data row (row 86) to the CompaNanny database sheet, matching
# the existing "plain data row" look (no explicit cell style) used by the
# other rows in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 86

$ws.Cells.Item($row, 1).Value = "Partou"
$ws.Cells.Item($row, 2).Value = "Partou KDV Blaricummerstraat 1A"
$ws.Cells.Item($row, 3).Value = "KDV"

# Column D holds a date-like value, but in this workbook it is stored as
# plain text (e.g. "2024-06-25"), not a real Excel date. Force the cell to
# text formatting first so Excel doesn't auto-convert the string into a
# date serial number, then restore the cell's style to "Normal" so no
# extra formatting/style index is left behind on the cell.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-08-27"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
